$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)
$ws.Range("H5").Value = "'2012-04-30"
$ws.Range("H5").ClearFormats()
$wb.Save()
